# Add a "Turkey" market test-data sheet, cloned from the existing "Spain"
# sheet (same layout/styles), then re-pointed at new shared-string values
# and tidied up so it becomes the new active tab.

$wb = $excel.ActiveWorkbook
$spain = $wb.Worksheets.Item("Spain")

# Duplicate "Spain" and place the copy immediately after it.
$spain.Copy($null, $spain) | Out-Null
$turkey = $wb.Worksheets.Item("Spain (2)")
$turkey.Name = "Turkey"

# Market name / Jira ticket for the new country.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3307"

# The copied rows 3-5 inherited Spain's taller (wrapped) row height;
# restore them to the sheet's normal auto height.
$turkey.Rows.Item(3).EntireRow.AutoFit() | Out-Null
$turkey.Rows.Item(4).EntireRow.AutoFit() | Out-Null
$turkey.Rows.Item(5).EntireRow.AutoFit() | Out-Null

# Widen column D a bit on the new sheet.
$turkey.Columns.Item(4).ColumnWidth = 25.6

# Put the cursor/selection on D7 for the new sheet (matches how it was
# left selected when saved).
$turkey.Range("D7").Select() | Out-Null

# "Spain" is no longer the active tab; clear its old single-cell
# selection in favor of selecting the whole used range.
$spain.Range("A1:D18").Select() | Out-Null

# Re-select Turkey last so it ends up the active/visible tab on reopen.
$turkey.Select() | Out-Null
$turkey.Range("D7").Select() | Out-Null
